$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Vessel Details" (sheet1 / Table1): add a new "Vessel Info Link" column
# ---------------------------------------------------------------------------
$wsVessel = $wb.Worksheets.Item("Vessel Details")
$loVessel = $wsVessel.ListObjects.Item("Table1")

# Add a new column to the table (will become column H)
$loVessel.ListColumns.Add() | Out-Null
$wsVessel.Range("G1").Copy($wsVessel.Range("H1")) | Out-Null
$wsVessel.Range("H1").Value = "Vessel Info Link"
$wsVessel.Range("H1").ColumnWidth = 109.26

# Populate the Vessel Info Link cells with hyperlinks
$wsVessel.Range("H2").Value = "https://www.marinetraffic.com/en/ais/details/ships/shipid:665124/mmsi:355104000/imo:9196319/vessel:SEA_CRUISER_2"
$wsVessel.Hyperlinks.Add($wsVessel.Range("H2"), "https://www.marinetraffic.com/en/ais/details/ships/shipid:665124/mmsi:355104000/imo:9196319/vessel:SEA_CRUISER_2") | Out-Null
$wsVessel.Range("H2").Style = "Hyperlink"

$wsVessel.Range("H3").Value = "https://www.marinetraffic.com/en/ais/details/ships/shipid:146379/mmsi:212882000/imo:9195133/vessel:EMS_HIGHWAY"
$wsVessel.Hyperlinks.Add($wsVessel.Range("H3"), "https://www.marinetraffic.com/en/ais/details/ships/shipid:146379/mmsi:212882000/imo:9195133/vessel:EMS_HIGHWAY") | Out-Null
$wsVessel.Range("H3").Style = "Hyperlink"

$wsVessel.Range("H4").Value = "https://www.marinetraffic.com/en/ais/details/ships/shipid:376403/mmsi:311996000/imo:9316294/vessel:THAMES_HIGHWAY"
$wsVessel.Hyperlinks.Add($wsVessel.Range("H4"), "https://www.marinetraffic.com/en/ais/details/ships/shipid:376403/mmsi:311996000/imo:9316294/vessel:THAMES_HIGHWAY") | Out-Null
$wsVessel.Range("H4").Style = "Hyperlink"

$wsVessel.Range("A1:H1").Select()

# ---------------------------------------------------------------------------
# Sheet "Deck Heights" (sheet2 / Table4): insert a new "Deck Id" column,
# rename "Deck Number" -> "Deck", and update a few deck labels
# ---------------------------------------------------------------------------
$wsDeck = $wb.Worksheets.Item("Deck Heights")

# Insert a blank column before the existing "Deck Number" column (column B)
$wsDeck.Columns("B:B").Insert()

$loDeck = $wsDeck.ListObjects.Item("Table4")
$loDeck.Resize($wsDeck.Range("A1:F25"))

# Restore / set all header names explicitly (keeps table column metadata in sync)
$wsDeck.Range("A1").Value = "Vessel"
$wsDeck.Range("B1").Value = "Deck Id"
$wsDeck.Range("C1").Value = "Deck"
$wsDeck.Range("D1").Value = "Average Deck Height (m)"
$wsDeck.Range("E1").Value = "Deck Type"
$wsDeck.Range("F1").Value = "Notes"

$wsDeck.Range("B1").ColumnWidth = 9.14

# Fill in the new "Deck Id" sequential numbers per vessel group
$deckIds = @(1,2,3,4,5, 1,2,3,4,5,6,7, 1,2,3,4,5,6,7,8,9,10,11,12)
for ($i = 0; $i -lt $deckIds.Length; $i++) {
    $row = $i + 2
    $wsDeck.Cells.Item($row, 2).Value = $deckIds[$i]
}

# Update specific "Deck" labels (Thames Highway rows) to match new values
$wsDeck.Range("C14").Value = "B02"
$wsDeck.Range("C15").Value = "B01"
$wsDeck.Range("C25").Value = "WD"

$wsDeck.Range("A1:F1").Select()

# ---------------------------------------------------------------------------
# Sheet "Comments" (sheet3): no data changes, just refresh the selection
# ---------------------------------------------------------------------------
$wsComments = $wb.Worksheets.Item("Comments")
$wsComments.Range("A1:E1").Select()

$wsVessel.Select()
